# Appends 5 new species-observation rows (23-27) to the "Artfynd" sheet,
# reproducing the upstream data export exactly (values + cell types).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 23 ---
$ws.Range("I23").NumberFormat = "@"
$ws.Range("Y23").NumberFormat = "@"
$ws.Range("AA23").NumberFormat = "@"
$ws.Range("A23").Value = 112501059
$ws.Range("B23").Value = 96735
$ws.Range("C23").Value = 'Ovaliderad'
$ws.Range("D23").Value = 'VU'
$ws.Range("E23").Value = 220787
$ws.Range("F23").Value = 'Knärot'
$ws.Range("G23").Value = 'Goodyera repens'
$ws.Range("H23").Value = '(L.) R. Br.'
$ws.Range("I23").Value = '40'
$ws.Range("J23").Value = 'stjälkar/strån/skott'
$ws.Range("K23").Value = 'fullt utvecklade blad'
$ws.Range("P23").Value = 'Nybygget S 645 m, Ög'
$ws.Range("Q23").Value = 562710
$ws.Range("R23").Value = 6504626
$ws.Range("S23").Value = 10
$ws.Range("T23").Value = 'Östergötland'
$ws.Range("U23").Value = 'Norrköping'
$ws.Range("V23").Value = 'Östergötland'
$ws.Range("W23").Value = 'Kvillinge'
$ws.Range("Y23").Value = '2023-10-02'
$ws.Range("AA23").Value = '2023-10-02'
$ws.Range("AD23").Value = $false
$ws.Range("AE23").Value = $false
$ws.Range("AG23").Value = $false
$ws.Range("AH23").Value = 'Barrskog'
$ws.Range("AW23").Value = 'Mirjam Ideström'
$ws.Range("AX23").Value = 'Mirjam Ideström'

# --- Row 24 ---
$ws.Range("I24").NumberFormat = "@"
$ws.Range("Y24").NumberFormat = "@"
$ws.Range("AA24").NumberFormat = "@"
$ws.Range("A24").Value = 112500915
$ws.Range("B24").Value = 103781
$ws.Range("C24").Value = 'Ovaliderad'
$ws.Range("D24").Value = 'LC'
$ws.Range("E24").Value = 221144
$ws.Range("F24").Value = 'Grönpyrola'
$ws.Range("G24").Value = 'Pyrola chlorantha'
$ws.Range("H24").Value = 'Sw.'
$ws.Range("I24").Value = '20'
$ws.Range("J24").Value = 'plantor/tuvor'
$ws.Range("K24").Value = 'fullt utvecklade blad'
$ws.Range("P24").Value = 'Nybygget S 767 m, Ög'
$ws.Range("Q24").Value = 562650
$ws.Range("R24").Value = 6504519
$ws.Range("S24").Value = 10
$ws.Range("T24").Value = 'Östergötland'
$ws.Range("U24").Value = 'Norrköping'
$ws.Range("V24").Value = 'Östergötland'
$ws.Range("W24").Value = 'Kvillinge'
$ws.Range("Y24").Value = '2023-10-02'
$ws.Range("AA24").Value = '2023-10-02'
$ws.Range("AD24").Value = $false
$ws.Range("AE24").Value = $false
$ws.Range("AG24").Value = $false
$ws.Range("AH24").Value = 'Barrskog'
$ws.Range("AW24").Value = 'Mirjam Ideström'
$ws.Range("AX24").Value = 'Mirjam Ideström'

# --- Row 25 ---
$ws.Range("I25").NumberFormat = "@"
$ws.Range("Y25").NumberFormat = "@"
$ws.Range("AA25").NumberFormat = "@"
$ws.Range("A25").Value = 112500988
$ws.Range("B25").Value = 103781
$ws.Range("C25").Value = 'Ovaliderad'
$ws.Range("D25").Value = 'LC'
$ws.Range("E25").Value = 221144
$ws.Range("F25").Value = 'Grönpyrola'
$ws.Range("G25").Value = 'Pyrola chlorantha'
$ws.Range("H25").Value = 'Sw.'
$ws.Range("I25").Value = '15'
$ws.Range("J25").Value = 'plantor/tuvor'
$ws.Range("K25").Value = 'i frukt'
$ws.Range("P25").Value = 'Nybygget S 712 m, Ög'
$ws.Range("Q25").Value = 562650
$ws.Range("R25").Value = 6504577
$ws.Range("S25").Value = 10
$ws.Range("T25").Value = 'Östergötland'
$ws.Range("U25").Value = 'Norrköping'
$ws.Range("V25").Value = 'Östergötland'
$ws.Range("W25").Value = 'Kvillinge'
$ws.Range("Y25").Value = '2023-10-02'
$ws.Range("AA25").Value = '2023-10-02'
$ws.Range("AD25").Value = $false
$ws.Range("AE25").Value = $false
$ws.Range("AG25").Value = $false
$ws.Range("AH25").Value = 'Barrskog'
$ws.Range("AW25").Value = 'Mirjam Ideström'
$ws.Range("AX25").Value = 'Mirjam Ideström'

# --- Row 26 ---
$ws.Range("Y26").NumberFormat = "@"
$ws.Range("AA26").NumberFormat = "@"
$ws.Range("A26").Value = 112500774
$ws.Range("B26").Value = 90806
$ws.Range("C26").Value = 'Ovaliderad'
$ws.Range("D26").Value = 'NT'
$ws.Range("E26").Value = 4361
$ws.Range("F26").Value = 'Orange taggsvamp'
$ws.Range("G26").Value = 'Hydnellum aurantiacum'
$ws.Range("H26").Value = '(Batsch:Fr.) P.Karst.'
$ws.Range("P26").Value = 'Lilla Bergsätter VNV 740 m, Ög'
$ws.Range("Q26").Value = 562914
$ws.Range("R26").Value = 6504341
$ws.Range("S26").Value = 10
$ws.Range("T26").Value = 'Östergötland'
$ws.Range("U26").Value = 'Norrköping'
$ws.Range("V26").Value = 'Östergötland'
$ws.Range("W26").Value = 'Kvillinge'
$ws.Range("Y26").Value = '2023-10-02'
$ws.Range("AA26").Value = '2023-10-02'
$ws.Range("AD26").Value = $false
$ws.Range("AE26").Value = $false
$ws.Range("AG26").Value = $false
$ws.Range("AH26").Value = 'Barrskog'
$ws.Range("AW26").Value = 'Mirjam Ideström'
$ws.Range("AX26").Value = 'Mirjam Ideström'

# --- Row 27 ---
$ws.Range("Y27").NumberFormat = "@"
$ws.Range("AA27").NumberFormat = "@"
$ws.Range("A27").Value = 112499465
$ws.Range("B27").Value = 90816
$ws.Range("C27").Value = 'Ovaliderad'
$ws.Range("D27").Value = 'NT'
$ws.Range("E27").Value = 788
$ws.Range("F27").Value = 'Gul taggsvamp'
$ws.Range("G27").Value = 'Hydnellum geogenium'
$ws.Range("H27").Value = '(Fr.) Banker'
$ws.Range("P27").Value = 'Skärlötamarken, Ög'
$ws.Range("Q27").Value = 562724
$ws.Range("R27").Value = 6504540
$ws.Range("S27").Value = 75
$ws.Range("T27").Value = 'Östergötland'
$ws.Range("U27").Value = 'Norrköping'
$ws.Range("V27").Value = 'Östergötland'
$ws.Range("W27").Value = 'Kvillinge'
$ws.Range("Y27").Value = '2023-09-29'
$ws.Range("AA27").Value = '2023-09-29'
$ws.Range("AC27").Value = 'I barrskog.'
$ws.Range("AD27").Value = $false
$ws.Range("AE27").Value = $false
$ws.Range("AG27").Value = $false
$ws.Range("AW27").Value = 'Jens Johannesson'
$ws.Range("AX27").Value = 'Jens Johannesson'

Write-Host "Added rows 23-27 to Artfynd sheet."
